$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data: date 2012-09-17 (serial 41169) with new activity text
$ws.Range("A23").Value = 41169
$ws.Range("B23").Value = "LocalScan with optimized memory access works now. Is somehow still slower =/"

# Update the selected cell to B29
[void]$ws.Range("B29").Select()
